$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Delete the row for account 004508516 (EDUARDO) first (higher row index),
# then the row for account 004805273 (CLISIA), so row indices remain valid.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()
